$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (preserve rich-text run formatting where possible) ---
$volCell = $ws.Range("A8")
$volCell.Characters(21, 2).Text = "17"

$weekCell = $ws.Range("C9")
$weekCell.Characters(27, 9).Text = "4/24/2023"
$weekCell.Characters(47, 9).Text = "4/30/2023"

# --- Weekly crime-stat table updates (rows 14-30) ---
# Row 14
$ws.Range("C14").Value = 2
$ws.Range("E14").Value = 100
$ws.Range("G14").Value = 7
$ws.Range("H14").Value = -28.571428571428
$ws.Range("I14").Value = 21
$ws.Range("J14").Value = 18
$ws.Range("K14").Value = 16.666666666666
$ws.Range("L14").Value = -16
$ws.Range("M14").Value = -46.153846153846
$ws.Range("N14").Value = -85.81081081081

# Row 15
$ws.Range("C15").Value = 6
$ws.Range("E15").Value = 100
$ws.Range("F15").Value = 18
$ws.Range("G15").Value = 12
$ws.Range("H15").Value = 50
$ws.Range("I15").Value = 78
$ws.Range("J15").Value = 75
$ws.Range("K15").Value = 4
$ws.Range("L15").Value = 9.859154929577
$ws.Range("M15").Value = 1.298701298701
$ws.Range("N15").Value = -59.162303664921

# Row 16
$ws.Range("C16").Value = 51
$ws.Range("D16").Value = 45
$ws.Range("E16").Value = 13.333333333333
$ws.Range("F16").Value = 181
$ws.Range("G16").Value = 172
$ws.Range("H16").Value = 5.232558139534
$ws.Range("I16").Value = 750
$ws.Range("J16").Value = 770
$ws.Range("K16").Value = -2.597402597402
$ws.Range("L16").Value = 29.757785467128
$ws.Range("M16").Value = -26.108374384236
$ws.Range("N16").Value = -85.738733599543

# Row 17
$ws.Range("C17").Value = 69
$ws.Range("D17").Value = 81
$ws.Range("E17").Value = -14.814814814814
$ws.Range("F17").Value = 314
$ws.Range("G17").Value = 306
$ws.Range("H17").Value = 2.614379084967
$ws.Range("I17").Value = 1257
$ws.Range("J17").Value = 1226
$ws.Range("K17").Value = 2.52854812398
$ws.Range("L17").Value = 28.923076923076
$ws.Range("M17").Value = 28.134556574923
$ws.Range("N17").Value = -50.492319810949

# Row 18
$ws.Range("C18").Value = 30
$ws.Range("D18").Value = 56
$ws.Range("E18").Value = -46.428571428571
$ws.Range("F18").Value = 141
$ws.Range("G18").Value = 174
$ws.Range("H18").Value = -18.965517241379
$ws.Range("I18").Value = 676
$ws.Range("J18").Value = 791
$ws.Range("K18").Value = -14.538558786346
$ws.Range("L18").Value = 11.735537190082
$ws.Range("M18").Value = -20.188902007083
$ws.Range("N18").Value = -81.759309228278

# Row 19
$ws.Range("C19").Value = 99
$ws.Range("D19").Value = 122
$ws.Range("E19").Value = -18.852459016393
$ws.Range("F19").Value = 434
$ws.Range("G19").Value = 418
$ws.Range("H19").Value = 3.827751196172
$ws.Range("I19").Value = 1832
$ws.Range("J19").Value = 1758
$ws.Range("K19").Value = 4.209328782707
$ws.Range("L19").Value = 40.814757878555
$ws.Range("M19").Value = 50.163934426229
$ws.Range("N19").Value = -8.170426065162

# Row 20
$ws.Range("C20").Value = 28
$ws.Range("D20").Value = 30
$ws.Range("E20").Value = -6.666666666666
$ws.Range("F20").Value = 134
$ws.Range("G20").Value = 118
$ws.Range("H20").Value = 13.559322033898
$ws.Range("I20").Value = 535
$ws.Range("J20").Value = 566
$ws.Range("K20").Value = -5.47703180212
$ws.Range("L20").Value = 38.242894056847
$ws.Range("M20").Value = 26.777251184834
$ws.Range("N20").Value = -83.492749151496

# Row 21
$ws.Range("C21").Value = 285
$ws.Range("D21").Value = 338
$ws.Range("E21").Value = -15.680473372781
$ws.Range("F21").Value = 1227
$ws.Range("G21").Value = 1207
$ws.Range("H21").Value = 1.6570008285
$ws.Range("I21").Value = 5149
$ws.Range("J21").Value = 5204
$ws.Range("K21").Value = -1.056879323597
$ws.Range("L21").Value = 30.618975139523
$ws.Range("M21").Value = 11.910454249076
$ws.Range("N21").Value = -69.851864863282

# Row 22
$ws.Range("C22").Value = 5
$ws.Range("D22").Value = 10
$ws.Range("E22").Value = -50
$ws.Range("F22").Value = 27
$ws.Range("G22").Value = 27
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 103
$ws.Range("J22").Value = 127
$ws.Range("K22").Value = -18.897637795275
$ws.Range("L22").Value = 33.766233766233
$ws.Range("M22").Value = -23.134328358209

# Row 23
$ws.Range("C23").Value = 25
$ws.Range("D23").Value = 29
$ws.Range("E23").Value = -13.793103448275
$ws.Range("F23").Value = 118
$ws.Range("G23").Value = 102
$ws.Range("H23").Value = 15.686274509803
$ws.Range("I23").Value = 521
$ws.Range("J23").Value = 458
$ws.Range("K23").Value = 13.755458515283
$ws.Range("L23").Value = 20.046082949308
$ws.Range("M23").Value = 68.064516129032

# Row 24
$ws.Range("C24").Value = 262
$ws.Range("D24").Value = 293
$ws.Range("E24").Value = -10.580204778157
$ws.Range("F24").Value = 942
$ws.Range("G24").Value = 1052
$ws.Range("H24").Value = -10.456273764258
$ws.Range("I24").Value = 3936
$ws.Range("J24").Value = 4025
$ws.Range("K24").Value = -2.211180124223
$ws.Range("L24").Value = 26.722472633612
$ws.Range("M24").Value = 29.943875866622

# Row 25
$ws.Range("C25").Value = 117
$ws.Range("D25").Value = 141
$ws.Range("E25").Value = -17.021276595744
$ws.Range("F25").Value = 447
$ws.Range("G25").Value = 473
$ws.Range("H25").Value = -5.496828752642
$ws.Range("I25").Value = 1875
$ws.Range("J25").Value = 1915
$ws.Range("K25").Value = -2.088772845953
$ws.Range("L25").Value = 45.123839009287
$ws.Range("M25").Value = -21.416596814752

# Row 26
$ws.Range("D26").Value = 5
$ws.Range("E26").Value = 20
$ws.Range("F26").Value = 21
$ws.Range("G26").Value = 17
$ws.Range("H26").Value = 23.529411764705
$ws.Range("I26").Value = 107
$ws.Range("J26").Value = 115
$ws.Range("K26").Value = -6.95652173913
$ws.Range("L26").Value = -12.295081967213

# Row 27
$ws.Range("C27").Value = 15
$ws.Range("D27").Value = 17
$ws.Range("E27").Value = -11.764705882352
$ws.Range("F27").Value = 55
$ws.Range("G27").Value = 49
$ws.Range("H27").Value = 12.244897959183
$ws.Range("I27").Value = 199
$ws.Range("J27").Value = 191
$ws.Range("K27").Value = 4.188481675392
$ws.Range("L27").Value = -4.326923076923

# Row 28
$ws.Range("C28").Value = 3
$ws.Range("E28").Value = -40
$ws.Range("F28").Value = 12
$ws.Range("G28").Value = 28
$ws.Range("H28").Value = -57.142857142857
$ws.Range("I28").Value = 70
$ws.Range("J28").Value = 79
$ws.Range("K28").Value = -11.392405063291
$ws.Range("L28").Value = -32.038834951456
$ws.Range("M28").Value = -50.704225352112
$ws.Range("N28").Value = -88.543371522094

# Row 29
$ws.Range("C29").Value = 3
$ws.Range("D29").Value = 4
$ws.Range("E29").Value = -25
$ws.Range("F29").Value = 11
$ws.Range("G29").Value = 23
$ws.Range("H29").Value = -52.173913043478
$ws.Range("I29").Value = 60
$ws.Range("J29").Value = 71
$ws.Range("K29").Value = -15.492957746478
$ws.Range("L29").Value = -33.333333333333
$ws.Range("M29").Value = -43.925233644859
$ws.Range("N29").Value = -89.266547406082

# Row 30
$ws.Range("C30").Value = 1
$ws.Range("E30").Value = 0
$ws.Range("G30").Value = 3
$ws.Range("H30").Value = 33.333333333333
$ws.Range("I30").Value = 28
$ws.Range("J30").Value = 23
$ws.Range("K30").Value = 21.739130434782
$ws.Range("L30").Value = 133.333333333333

